$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 259
$ws.Range("I5").Value = 112.8
$ws.Range("K5").Value = 112.8
$ws.Range("M5").Value = 2.200000000000003
$ws.Range("H17").Value = 836.63635
$ws.Range("J17").Value = 800.4286
$ws.Range("L17").Value = 2401.2858
$ws.Range("N17").Value = -2737.2858
$ws.Range("H28").Value = 1675.75
$ws.Range("I28").Value = 1932.3334
$ws.Range("J28").Value = 906
$ws.Range("K28").Value = 1932.3334
$ws.Range("L28").Value = 906
$ws.Range("M28").Value = -1447.3334
$ws.Range("N28").Value = -1876
$ws.Range("H86").Value = 5585
$ws.Range("I86").Value = 4383.1665
$ws.Range("J86").Value = 7387.75
$ws.Range("K86").Value = 4383.1665
$ws.Range("L86").Value = 7387.75
$ws.Range("M86").Value = -3260.1665
$ws.Range("N86").Value = -9633.75
$ws.Range("H89").Value = 5585
$ws.Range("I89").Value = 4383.1665
$ws.Range("J89").Value = 7387.75
$ws.Range("K89").Value = 21915.8325
$ws.Range("L89").Value = 36938.75
$ws.Range("M89").Value = -16299.8325
$ws.Range("N89").Value = -48170.75
$ws.Range("H106").Value = 2300
$ws.Range("I106").Value = 2300
$ws.Range("K106").Value = 2300
$ws.Range("M106").Value = -1669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 7300
$ws.Range("J27").Value = 7300
$ws.Range("L27").Value = 7300
$ws.Range("N27").Value = -7668
$ws.Range("H45").Value = 3873.7727
$ws.Range("I45").Value = 1278.3334
$ws.Range("K45").Value = 1278.3334
$ws.Range("M45").Value = -901.3334
$ws.Range("H61").Value = 1993.4546
$ws.Range("I61").Value = 1993.4546
$ws.Range("K61").Value = 1993.4546
$ws.Range("M61").Value = -1781.4546
$ws.Range("H88").Value = 2733.7
$ws.Range("J88").Value = 3034.8572
$ws.Range("L88").Value = 3034.8572
$ws.Range("N88").Value = -3846.8572
$ws.Range("H91").Value = 2733.7
$ws.Range("J91").Value = 3034.8572
$ws.Range("L91").Value = 3034.8572
$ws.Range("N91").Value = -5842.8572
$ws.Range("H102").Value = 1435.5714
$ws.Range("I102").Value = 1435.5714
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1435.5714
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 186.4286
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 608.36365
$ws.Range("I122").Value = 628.8
$ws.Range("J122").Value = 404
$ws.Range("K122").Value = 1886.4
$ws.Range("L122").Value = 1212
$ws.Range("M122").Value = 563.6000000000001
$ws.Range("N122").Value = -6112
$ws.Range("H132").Value = 4506
$ws.Range("I132").Value = 4506
$ws.Range("K132").Value = 13518
$ws.Range("M132").Value = -10988
$ws.Range("H136").Value = 1993.4546
$ws.Range("I136").Value = 1993.4546
$ws.Range("K136").Value = 5980.3638
$ws.Range("M136").Value = -3430.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 950
$ws.Range("I20").Value = 950
$ws.Range("K20").Value = 950
$ws.Range("M20").Value = -703
$ws.Range("H54").Value = 11000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H94").Value = 3554
$ws.Range("I94").Value = 3006.5386
$ws.Range("K94").Value = 3006.5386
$ws.Range("M94").Value = -2555.5386
$ws.Range("H134").Value = 2681.6428
$ws.Range("I134").Value = 2681.6428
$ws.Range("K134").Value = 8044.928400000001
$ws.Range("M134").Value = -5509.928400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 449.46155
$ws.Range("I22").Value = 477
$ws.Range("K22").Value = 477
$ws.Range("M22").Value = -127
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2500
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 12500
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -18740
$ws.Range("H99").Value = 4999.9
$ws.Range("I99").Value = 4999.9
$ws.Range("K99").Value = 4999.9
$ws.Range("M99").Value = -3501.9
$ws.Range("H126").Value = 4999.9
$ws.Range("I126").Value = 4999.9
$ws.Range("K126").Value = 14999.7
$ws.Range("M126").Value = -12529.7
$ws.Range("H132").Value = 2447.1304
$ws.Range("I132").Value = 2376.5454
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7129.6362
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4599.6362
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 884.2857
$ws.Range("I50").Value = 350
$ws.Range("J50").Value = 973.3333
$ws.Range("K50").Value = 1050
$ws.Range("L50").Value = 2919.9999
$ws.Range("M50").Value = -569
$ws.Range("N50").Value = -3881.9999
$ws.Range("H53").Value = 884.2857
$ws.Range("I53").Value = 350
$ws.Range("J53").Value = 973.3333
$ws.Range("K53").Value = 1050
$ws.Range("L53").Value = 2919.9999
$ws.Range("M53").Value = -569
$ws.Range("N53").Value = -3881.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8506
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 8506
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H122").Value = 1538.7778
$ws.Range("I122").Value = 1593.75
$ws.Range("J122").Value = 1099
$ws.Range("K122").Value = 4781.25
$ws.Range("L122").Value = 3297
$ws.Range("M122").Value = -2331.25
$ws.Range("N122").Value = -8197
$ws.Range("H132").Value = 4433.25
$ws.Range("I132").Value = 5462.375
$ws.Range("K132").Value = 16387.125
$ws.Range("M132").Value = -13857.125
$ws.Range("H134").Value = 26000
$ws.Range("J134").Value = 26000
$ws.Range("L134").Value = 78000
$ws.Range("N134").Value = -83070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12040.75
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 12040.75
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 12040.75
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -12630.75
$ws.Range("H27").Value = 12040.75
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 12040.75
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 12040.75
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -12254.75
$ws.Range("H82").Value = 478.6
$ws.Range("I82").Value = 397
$ws.Range("J82").Value = 533
$ws.Range("K82").Value = 397
$ws.Range("L82").Value = 533
$ws.Range("M82").Value = -36
$ws.Range("N82").Value = -1255
$ws.Range("H85").Value = 478.6
$ws.Range("I85").Value = 397
$ws.Range("J85").Value = 533
$ws.Range("K85").Value = 397
$ws.Range("L85").Value = 533
$ws.Range("M85").Value = 851
$ws.Range("N85").Value = -3029
$ws.Range("H132").Value = 6969
$ws.Range("J132").Value = 6969
$ws.Range("L132").Value = 20907
$ws.Range("N132").Value = -25967

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 5000
$ws.Range("I47").Value = 5000
$ws.Range("K47").Value = 5000
$ws.Range("M47").Value = -4428
$ws.Range("H81").Value = 1484.5714
$ws.Range("I81").Value = 1484.5714
$ws.Range("K81").Value = 2969.1428
$ws.Range("M81").Value = -1908.1428
$ws.Range("H84").Value = 1484.5714
$ws.Range("I84").Value = 1484.5714
$ws.Range("K84").Value = 14845.714
$ws.Range("M84").Value = -9541.714
$ws.Range("H126").Value = 2735.875
$ws.Range("J126").Value = 5201
$ws.Range("L126").Value = 15603
$ws.Range("N126").Value = -20543
